$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.028.80"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "3.515.28"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'605.09"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "'147.87"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("D7").Value = "3.516.33"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Value = "'7.83"
$ws.Range("E11").Value = "  +3.13%  "
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "4.109.80"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "'31.60"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").Value = "3.512.62"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "67.006.28"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'10.65"
$ws.Range("E19").Value = "  +8.06%  "
$ws.Range("D20").Value = "'6.40"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "'15.38"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'435.88"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "'0.611"
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").Value = "'79.58"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("D25").Value = "3.653.11"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -4.01%  "
$ws.Range("D28").Value = "'9.86"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").Value = "'8.33"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'1.59"
$ws.Range("E31").Value = "  -3.38%  "
$ws.Range("D32").Value = "'0.168"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'25.38"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "3.508.84"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").Value = "'5.90"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "'8.02"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'0.0892"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "'169.44"
$ws.Range("E42").Value = "  -3.25%  "
$ws.Range("E43").Value = "  -9.52%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'0.896"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").Value = "'28.83"
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("D47").Value = "'45.76"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'1.32"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").Value = "'7.48"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").Value = "'0.989"
$ws.Range("E51").Value = "  -0.44%  "
